# Revert "Unidad 2 actualizado 19 jun 2013"
$p = $ppt.ActivePresentation

# 1. Remove the last two slides (slide 9 "Ciclo while" and slide 10 "Ciclo do-while")
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()

# 2. Roll back the cached date field text on the slide master and every slide layout
#    from 6/19/2013 back to 6/18/2013.
$p.SlideMaster.Shapes.Item(4).TextFrame.TextRange.Text = "6/18/2013"

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "6/19/2013") {
                $shape.TextFrame.TextRange.Text = "6/18/2013"
            }
        }
    }
}

# 3. Split the "UNIDAD II" run on the title slide into two runs: "UNIDAD " and "II"
$s1 = $p.Slides.Item(1)
for ($k = 1; $k -le $s1.Shapes.Count; $k++) {
    $shape = $s1.Shapes.Item($k)
    if ($shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        $idx = $fullText.IndexOf("UNIDAD II")
        if ($idx -ge 0) {
            # 1-based Characters(Start, Length); re-assigning just the "UNIDAD "
            # prefix (not the whole run) forces the engine to split the run into
            # "UNIDAD " + "II" instead of leaving one merged run.
            $part = $tr.Characters($idx + 1, 7)
            $part.Text = "UNIDAD "
        }
    }
}
